$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 16.45659999999999
$ws.Range("E12").Value = 17.9906
$ws.Range("E18").Value = 17.66510000000002
$ws.Range("E37").Value = 16.68320000000001
$ws.Range("E55").Value = 16.64460000000001
$ws.Range("E68").Value = 17.21150000000002
$ws.Range("E77").Value = 18.21230000000003
$ws.Range("E78").Value = 16.51870000000003
